$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly records go in at the top (row 2/3); the existing 17 data
# rows (old rows 2-18) shift down by two (to rows 4-20).
$ws.Rows("2:3").Insert()

# The Insert() above drags the bold header formatting down onto the new
# rows - clear that so the new rows look like ordinary data rows again.
$ws.Range("A2:T3").ClearFormats()

# Column D (Fecha) keeps the same date display format used by every other
# data row in the sheet.
$ws.Range("D2:D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 2 - new Damasco / Castle Brite / Primera record for Región de O'Higgins
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "Macroferia Regional de Talca"
$ws.Range("C2").Value = "Maule"
$ws.Range("D2").Value = 44530
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100103
$ws.Range("H2").Value = "Frutos de hueso (carozo)"
$ws.Range("I2").Value = 100103003
$ws.Range("J2").Value = "Damasco"
$ws.Range("K2").Value = "Castle Brite"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 130
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 20000
$ws.Range("Q2").Value = "`$/caja 15 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1333
$ws.Range("T2").Value = 15

# Row 3 - new Damasco / Castle Brite / Segunda record for Región de O'Higgins
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "Macroferia Regional de Talca"
$ws.Range("C3").Value = "Maule"
$ws.Range("D3").Value = 44530
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100103
$ws.Range("H3").Value = "Frutos de hueso (carozo)"
$ws.Range("I3").Value = 100103003
$ws.Range("J3").Value = "Damasco"
$ws.Range("K3").Value = "Castle Brite"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 150
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "`$/caja 15 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 15
